# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) inside specific bullet
# paragraphs of the resume, matching the target diff.

$d = $word.ActiveDocument

# Highlight color used across all formats: #2C3E50 -> Word BGR-packed
# integer (R + G*256 + B*65536) for Font.Color.
$highlightColor = 5258796

function Set-MetricHighlight {
    param($Paragraph, $Metric)

    $range = $Paragraph.Range
    $find = $range.Find
    $find.ClearFormatting()
    $find.Text = $Metric
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.Execute() | Out-Null

    if ($find.Found) {
        $range.Font.Bold = $true
        $range.Font.Color = $highlightColor
    }
}

# Locate target paragraphs by their (unique) original plain text, then
# bold+color the embedded metric substrings -- Word automatically splits
# the existing run(s) around the matched text when formatting is applied
# to a sub-range, producing the same run structure as the diff.

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text

    if ($text -like "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%*") {
        Set-MetricHighlight $p "23%"
        Set-MetricHighlight $p "64%"
    }
    elseif ($text -like "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%*") {
        Set-MetricHighlight $p "87%"
        Set-MetricHighlight $p "71%"
        Set-MetricHighlight $p "±4.2%"
        Set-MetricHighlight $p "±2.1%"
    }
    elseif ($text -like "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development*") {
        Set-MetricHighlight $p "1,200"
    }
    elseif ($text -like "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+*") {
        Set-MetricHighlight $p "`$400M"
        Set-MetricHighlight $p "`$1B"
    }
    elseif ($text -like "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M*") {
        Set-MetricHighlight $p "73.5%"
        Set-MetricHighlight $p "`$4.7M"
    }
    elseif ($text -like "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%*") {
        Set-MetricHighlight $p "87%"
        Set-MetricHighlight $p "71%"
    }
}

Write-Host "Metric highlighting applied."
